# Dataframe ST.xlsx - daily refresh edit
#
# 1. Rename product "LAYS CLASICAS 94GRX25" -> "LAYS CLASICAS 85GX25"
#    everywhere it occurs (Sheet1!B3, Sheet3!B3, Sheet3!A29).
# 2. Refresh the lookup table on Sheet3 (A20:B36) with the new metric
#    values for 24-nov (this ripples into Sheet1 columns CB/CC via the
#    existing VLOOKUP formulas).
# 3. Append a new "24-nov" snapshot column (CP) on Sheet1, hard-coded
#    with the same values the VLOOKUP now returns for each row.
# 4. Update the sheet's active selection to CP5, matching the saved
#    workbook state.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

$oldName = "LAYS CLASICAS 94GRX25"
$newName = "LAYS CLASICAS 85GX25"

# --- 1. Rename the product everywhere it appears ------------------------
$ws1.Range("B3").Value2 = $newName
$ws3.Range("B3").Value2 = $newName
$ws3.Range("A29").Value2 = $newName

# --- 2. Refresh Sheet3's lookup table (A20:B36) --------------------------
$lookup = @{
    "3D QUESO 92GX27"                     = 9.9712534292306678
    "CHEETOS QUESO 85GX24X1"              = 3.6639934842270918
    "DORITOS QUESO 129GX19"               = 8.4008654009124193
    "DORITOS QUESO 70X40G"                = 9.698103896120438
    "DORITOS QUESO 77GX26"                = 7.400326496359372
    "LAYS CEBOLLA CARAMELIZADA 85GX25"    = 0
    "LAYS CLASICAS 145GRX18"              = 16.602094538445346
    "LAYS CLASICAS 249GRX14"              = 6.3920077424003381
    "LAYS CLASICAS 40GX68"                = 6.5757244399532953
    "LAYS CLASICAS 85GX25"                = 11.634528865541975
    "LAYS ONDAS FH 30GX72"                = 7.5600000000233161
    "LAYS ONDAS FH 70GX28"                = 8.2979050000691572
    "LAYS QSO Y CEBOLLA 34GX72"           = 22.297874914106764
    "PEHUAMAR ACANALADA 520GX9"           = 6.9167179120300375
    "PEHUAMAR MAICITOS 285GX10"           = 11.438301023874551
    "PEHUAMAR PAPA LISA 520GX9"           = 10.164766619324679
    "QUAKER AVENA INSTANT FORTIF 18X280G" = 36.338314223202495
}

for ($r = 20; $r -le 36; $r++) {
    $prod = $ws3.Cells.Item($r, 1).Value2
    $ws3.Cells.Item($r, 2).Value2 = $lookup[$prod]
}

# --- 3. Append the new "24-nov" snapshot column (CP) on Sheet1 ----------
$ws1.Range("CP1").Value2 = "24-nov"
$ws1.Range("CP1").NumberFormat = "@"

$cpValues = @{
    2  = 6.5757244399532953
    3  = 11.634528865541975
    4  = 16.602094538445346
    5  = 6.3920077424003381
    6  = 9.698103896120438
    7  = 7.400326496359372
    8  = 8.4008654009124193
    9  = 10.164766619324679
    10 = 6.9167179120300375
    11 = 11.438301023874551
    12 = 9.9712534292306678
    13 = 3.6639934842270918
    14 = 36.338314223202495
    15 = 22.297874914106764
    16 = 0
    17 = 7.5600000000233161
    18 = 8.2979050000691572
}

for ($r = 2; $r -le 18; $r++) {
    $ws1.Cells.Item($r, 94).Value2 = $cpValues[$r]
    $ws1.Cells.Item($r, 94).NumberFormat = "0"
}

# --- 4. Restore the saved selection --------------------------------------
[void]$ws1.Activate()
[void]$ws1.Range("CP5").Select()
